$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at X:Y (everything from X onward shifts right by 2,
# i.e. old X/Y -> Z/AA, old AF -> AH, etc.)
$ws.Range("X1:Y1").EntireColumn.Insert()

# Row 9 (MHSTRF...) sdtm_annotation cell gets a "[NOT SUBMITTED];" prefix and the
# row grows taller to fit the longer wrapped text
$ws.Range("AH9").Value = "[NOT SUBMITTED];MHSTRF;MHSTRTPT;MHSTTPT"
$ws.Rows(9).RowHeight = 48

# Row 19 (MHENRF...) sdtm_annotation cell gets the same treatment (row height
# was already 48, so nothing else changes there)
$ws.Range("AH19").Value = "[NOT SUBMITTED];MHENRF;MHENRTPT;MHENTPT"

# New header cells for the inserted columns
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

Write-Output "done"
